$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: product_unit column is now right-aligned (predetermined discount
# column formatting pass touched the unit header too)
$ws.Range("F1").HorizontalAlignment = -4152

# Row 7 - product 15006 (Cadena plastica negra 10MM 3/8 X metro)
$ws.Range("B7").Value = 15006
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "Cadena plastica negra 10MM 3/8 X metro"
$ws.Range("E7").Value = 1070.8
$ws.Range("F7").Value = "Metro"
$ws.Range("F7").HorizontalAlignment = -4152
$ws.Range("G7").Value = "t"
$ws.Range("H7").Value = 13
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 15
$ws.Range("L7").Value = 5

# Row 8 - product 15007 (Lamina Prisma ACS de 3/4 #13MR (1,22X1,44))
$ws.Range("B8").Value = 15007
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = "Lamina Prisma ACS de 3/4 #13MR (1,22X1,44)"
$ws.Range("E8").Value = 20115.04
$ws.Range("F8").Value = "Unidad"
$ws.Range("F8").HorizontalAlignment = -4152
$ws.Range("G8").Value = "t"
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 15
$ws.Range("L8").Value = 5

# Row 9 - product 15008 (Polyacril)
$ws.Range("B9").Value = 15008
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = "Polyacril "
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "Unidad"
$ws.Range("F9").HorizontalAlignment = -4152
$ws.Range("G9").Value = "T"
$ws.Range("H9").Value = 13
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 15
